# Add season record columns (Wins, Losses, Ties) to the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (AC1) onto the
# three new header cells so they pick up the same bold/centered/bordered
# formatting (style index) used by the rest of row 1.
$ws.Range("AC1").Copy($ws.Range("AD1:AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every player row (2-52).
For ($row = 2; $row -le 52; $row++) {
    $ws.Cells.Item($row, 30).Value = 86
    $ws.Cells.Item($row, 31).Value = 76
    $ws.Cells.Item($row, 32).Value = 0
}
